# Clean up the NSW Population by Age Group dataset:
#  - Drop the old "ID" column (A), shifting Measure Name -> A and the
#    date/value columns left by one.
#  - Clear the stray trailing value left in the old C5/now-B5 cell.
#  - Add a new "NSW - Population Total" row with its total value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A entirely; B->A and C->B shift automatically.
$ws.Columns("A").Delete()

# The old total-row value (previously C5, now B5) is no longer wanted.
$ws.Range("B5").ClearContents()

# Append the new total row.
$ws.Range("A6").Value = "NSW - Population Total"
$ws.Range("B6").Value = 8167532

$ws.Range("B6").Select()
